$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @{ col letter = new value (as text) }
$updates = @{
    2 = @{ "D"="313.47"; "E"="2.34%"; "G"="20" }
    3 = @{ "D"="37.71"; "E"="1.11%"; "G"="20" }
    4 = @{ "D"="5.129"; "E"="0.18%"; "G"="20" }
    5 = @{ "D"="0.07914"; "E"="1.95%"; "G"="20" }
    6 = @{ "D"="4.426"; "E"="0.60%"; "G"="20" }
    7 = @{ "D"="1.927"; "E"="1.27%"; "G"="20" }
    8 = @{ "D"="8.288"; "E"="0.99%"; "G"="20" }
    9 = @{ "D"="2.862"; "E"="-9.79%"; "G"="20" }
    10 = @{ "D"="0.9223"; "E"="0.42%"; "G"="20" }
    11 = @{ "D"="0.1236"; "E"="-1.42%"; "G"="20" }
    12 = @{ "D"="0.1930"; "E"="2.43%"; "G"="20" }
    13 = @{ "D"="0.09266"; "E"="6.17%"; "G"="20" }
    14 = @{ "E"="-2.92%"; "G"="20" }
    15 = @{ "D"="0.09628"; "G"="20" }
    16 = @{ "D"="0.001384"; "E"="1.10%"; "G"="20" }
    17 = @{ "D"="0.005718"; "E"="-2.63%"; "G"="20" }
    18 = @{ "D"="3.513"; "E"="-0.99%"; "G"="20" }
    19 = @{ "E"="2.14%"; "G"="20" }
    20 = @{ "D"="5.265"; "E"="4.73%"; "G"="20" }
    21 = @{ "D"="0.1273"; "E"="-0.88%"; "G"="20" }
    22 = @{ "E"="3.64%"; "G"="20" }
    23 = @{ "E"="-0.48%"; "G"="20" }
    24 = @{ "D"="0.04363"; "E"="0.60%"; "G"="20" }
    25 = @{ "D"="0.001248"; "E"="2.21%"; "G"="20" }
    26 = @{ "D"="0.004318"; "E"="-3.81%"; "G"="20" }
    27 = @{ "D"="0.0001219"; "E"="-10.02%"; "G"="20" }
    28 = @{ "G"="20" }
    29 = @{ "G"="20" }
    30 = @{ "G"="20" }
    31 = @{ "G"="20" }
    32 = @{ "G"="20" }
    33 = @{ "G"="20" }
    34 = @{ "G"="20" }
    35 = @{ "G"="20" }
    36 = @{ "G"="20" }
    37 = @{ "G"="20" }
    38 = @{ "G"="20" }
    39 = @{ "D"="0.02234"; "E"="1.97%"; "G"="20" }
    40 = @{ "D"="0.05125"; "E"="3.99%"; "G"="20" }
    41 = @{ "D"="0.007446"; "E"="-3.28%"; "G"="20" }
    42 = @{ "E"="2.47%"; "G"="20" }
    43 = @{ "D"="0.008821"; "E"="-10.13%"; "G"="20" }
    44 = @{ "D"="0.001920"; "E"="-7.16%"; "G"="20" }
    45 = @{ "D"="0.008607"; "E"="-1.96%"; "G"="20" }
    46 = @{ "D"="0.00006729"; "E"="-1.44%"; "G"="20" }
    47 = @{ "E"="-0.43%"; "G"="20" }
    48 = @{ "D"="0.003345"; "E"="11.03%"; "G"="20" }
    49 = @{ "E"="-8.09%"; "G"="20" }
    50 = @{ "E"="-0.43%"; "G"="20" }
    51 = @{ "E"="-0.43%"; "G"="20" }
}

foreach ($row in $updates.Keys) {
    $cellUpdates = $updates[$row]
    foreach ($col in $cellUpdates.Keys) {
        $ws.Range("$col$row").Value = "'" + $cellUpdates[$col]
    }
}

Write-Host "Updated $($updates.Count) rows"